# S50001.MES.BIN translation sheet update
# Updates the "Edited" (F) column with new proposed translations / TL notes,
# the "Notes" (H) column with new translator comments, and corrects a
# couple of "Initial" (G) column lines that had typos/wording tweaks.
#
# Column layout: A=Status B=Block C=SpeakerID D=SpeakerName E=Japanese
#                F=Edited  G=Initial H=Notes   I/J=Other

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 'Maybe it''s okay to stop here for today...'
$ws.Range("H2").Value = 'Drop the ''already'' because it sounds very awkward.'
$ws.Range("F4").Value = 'Mana-chan happily begins to tidy up her writing utensils.'
$ws.Range("F5").Value = 'Although I feel unsatisfied, I realize I can''t do anything more after having been shown such perfect academic abilities.'
$ws.Range("H6").Value = 'Original text used "menu"'
$ws.Range("F10").Value = 'It''s not from a CD...'
$ws.Range("F11").Value = 'It''s FM...?'
$ws.Range("F12").Value = 'I start to relax a bit to the music, which is perfect for the afternoon.'
$ws.Range("F13").Value = 'Well, maybe it''s okay to enjoy the holiday feeling today...'
$ws.Range("F14").Value = 'Holiday...'
$ws.Range("G18").Value = 'Hm-?'
$ws.Range("F21").Value = 'It''s a voluntary day off...'
$ws.Range("G22").Value = '...Even so, you''re a student preparing for exams...'
$ws.Range("F23").Value = 'You''re so noisy. It''s fine as long as I can study properly, right?'
$ws.Range("F24").Value = 'The test earlier, wasn''t that the point of it? And since I was able to do it properly, there''s no reason to complain, right?'
$ws.Range("F25").Value = 'There''s no reason to complain, but...'
$ws.Range("F26").Value = '...What is it-?'
$ws.Range("F27").Value = 'N-never mind...'
$ws.Range("F29").Value = 'It''s okay. I''ll go shopping or something a while later.'
$ws.Range("F30").Value = '-I don''t want to be cooped up inside the house all day.'
$ws.Range("F31").Value = 'That''s not what this is about.'
$ws.Range("F32").Value = '...So, your mother or anyone else doesn''t say anything?'
$ws.Range("F33").Value = 'She doesn''t say anything, not to me...'
$ws.Range("F35").Value = '...Be quiet for a bit, Touya-san. ......I can''t hear the music...'
$ws.Range("F37").Value = 'I may have been a bit too intrusive.'
$ws.Range("F38").Value = 'I obediently say quiet and listen to the radio.'
$ws.Range("F44").Value = 'Well, sort of...'
$ws.Range("F45").Value = 'I give a vague answer.'
$ws.Range("F47").Value = 'After a short while, Mana-chan says,'
$ws.Range("G48").Value = 'Yeah. I wonder if it''s okay.'
